$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain plain text (many values look numeric,
# e.g. "223.61", and Excel would otherwise auto-convert them to numbers).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "27.308.22"
$ws.Range("E2").Value = "  -2.60%  "
$ws.Range("D3").Value = "1.707.08"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "223.61"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").Value = "0.5309"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "0.2654"
$ws.Range("E8").Value = "  -4.15%  "
$ws.Range("D9").Value = "0.06582"
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("D10").Value = "20.81"
$ws.Range("E10").Value = "  -4.08%  "
$ws.Range("D11").Value = "0.07624"
$ws.Range("E11").Value = "  -1.99%  "
$ws.Range("D12").Value = "4.566"
$ws.Range("E12").Value = "  -2.91%  "
$ws.Range("D13").Value = "1.700.58"
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").Value = "1.943.39"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "0.5718"
$ws.Range("E15").Value = "  -4.52%  "
$ws.Range("D16").Value = "0.0₅8159"
$ws.Range("E16").Value = "  -2.96%  "
$ws.Range("D17").Value = "67.70"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "27.319.54"
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("D19").Value = "215.60"
$ws.Range("E19").Value = "  -4.10%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "4.663"
$ws.Range("E21").Value = "  -3.70%  "
$ws.Range("D22").Value = "10.40"
$ws.Range("E22").Value = "  -4.54%  "
$ws.Range("D23").Value = "5.960"
$ws.Range("E23").Value = "  -4.46%  "
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "1.764"
$ws.Range("E25").Value = "  +7.19%  "
$ws.Range("D26").Value = "141.58"
$ws.Range("E26").Value = "  -3.16%  "
$ws.Range("D27").Value = "0.1213"
$ws.Range("E27").Value = "  -2.96%  "
$ws.Range("D28").Value = "7.263"
$ws.Range("E28").Value = "  -2.82%  "
$ws.Range("E29").Value = "  -4.97%  "
$ws.Range("D30").Value = "0.05399"
$ws.Range("E30").Value = "  -5.16%  "
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("D32").Value = "3.499"
$ws.Range("E32").Value = "  -6.00%  "
$ws.Range("D33").Value = "3.424"
$ws.Range("E33").Value = "  -2.98%  "
$ws.Range("D34").Value = "1.642"
$ws.Range("E34").Value = "  -2.04%  "
$ws.Range("D35").Value = "2.872"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").Value = "2.420"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").Value = "0.9477"
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("D38").Value = "0.5855"
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").Value = "0.01628"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("D40").Value = "5.865"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").Value = "1.045.39"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").Value = "1.003"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").Value = "0.8413"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("D44").Value = "100.71"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").Value = "1.850.07"
$ws.Range("E45").Value = "  -1.93%  "
$ws.Range("E46").Value = "  -1.43%  "
$ws.Range("D47").Value = "57.88"
$ws.Range("E47").Value = "  -3.62%  "
$ws.Range("D48").Value = "0.4503"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").Value = "1.005"
$ws.Range("E49").Value = "  +0.39%  "
$ws.Range("D50").Value = "8.074"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("E51").Value = "  -1.39%  "

# Restore the default (unstyled) appearance of the Price column now that
# the values are stored as text, so no style index is left on the cells.
$priceRange.Style = "Normal"
